$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The dataset gained two new HKL reflection entries ("Holden" and
# "Rizzie Spiral") which are inserted as rows 4 and 5 (right after the
# "Spiral5" row). All the previously existing data rows (old rows 4-29)
# shift down by two rows (to new rows 6-31), and the whole simulation was
# rerun, producing new values throughout. One existing label was also
# renamed ("Thomas Hex" -> "Matthies Hex").
# ---------------------------------------------------------------------------

# Step 1: shift existing data rows 4..29 down to 6..31 (process bottom-up so
# that we never overwrite a source row before it has been read).
for ($r = 29; $r -ge 4; $r--) {
    $dst = $r + 2
    for ($c = 1; $c -le 20; $c++) {
        $srcCell = $ws.Cells.Item($r, $c)
        $dstCell = $ws.Cells.Item($dst, $c)
        $dstCell.Value = $srcCell.Value2
    }
}

# Step 2: rows 30 and 31 are brand new rows beyond the old sheet bounds, so
# they do not inherit the bordered/bold style used by column A. Copy that
# formatting over from an existing styled cell in column A.
$ws.Range("A3").Copy() | Out-Null
$ws.Range("A30:A31").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Step 3: populate the two freshly-inserted rows (4 and 5) with the newly
# simulated data.
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "Holden"
$ws.Range("C4").Value = 0.9856866494859337
$ws.Range("D4").Value = 1.003172177542477
$ws.Range("E4").Value = 1.00380398576862
$ws.Range("F4").Value = 0.9856866494859337
$ws.Range("G4").Value = 1.003619783021732
$ws.Range("H4").Value = 1.003172177542477
$ws.Range("I4").Value = 0.9919814398771385
$ws.Range("J4").Value = 1.010264297135257
$ws.Range("K4").Value = 1.003172177542477
$ws.Range("L4").Value = 1.00380398576862
$ws.Range("M4").Value = 0.9947453176272767
$ws.Range("N4").Value = 0.9947453176272767
$ws.Range("O4").Value = 0.9938240250438973
$ws.Range("P4").Value = 0.9975542709323434
$ws.Range("Q4").Value = 0.9975542709323436
$ws.Range("R4").Value = 0.9989587475848769
$ws.Range("S4").Value = 0.9989587475848769
$ws.Range("T4").Value = 0.9997547221385261

$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "Rizzie Spiral"
$ws.Range("C5").Value = 0.8930423090184135
$ws.Range("D5").Value = 1.031078963800421
$ws.Range("E5").Value = 1.024328573234077
$ws.Range("F5").Value = 0.8930423090184135
$ws.Range("G5").Value = 1.026296604342213
$ws.Range("H5").Value = 1.031078963800421
$ws.Range("I5").Value = 0.9427355039141989
$ws.Range("J5").Value = 1.063590407046105
$ws.Range("K5").Value = 1.031078963800421
$ws.Range("L5").Value = 1.024328573234077
$ws.Range("M5").Value = 0.9586854411262455
$ws.Range("N5").Value = 0.9586854411262455
$ws.Range("O5").Value = 0.9533687953888966
$ws.Range("P5").Value = 0.9828166153509706
$ws.Range("Q5").Value = 0.9828166153509704
$ws.Range("R5").Value = 0.994882202463333
$ws.Range("S5").Value = 0.994882202463333
$ws.Range("T5").Value = 0.9968453935592381

# Step 4: rename "Thomas Hex" to "Matthies Hex" (this label now lives in
# row 11, since its old row 9 shifted down by two).
$ws.Range("B11").Value = "Matthies Hex"

# Step 5: the worksheet's used range/dimension grew by two rows (29 -> 31).
$ws.Range("A1:T31").Select() | Out-Null
